$wb = $excel.ActiveWorkbook
$dbd = $wb.Worksheets.Item("DBD")
$dbs = $wb.Worksheets.Item("DBS")

# Switch to the DBS sheet and add a new lookup-definition row (row 4),
# copying row 3's formatting down first, then filling in the new values.
$dbs.Select()

$dbs.Range("A3").Copy()
$dbs.Range("A4").PasteSpecial(-4122)
$dbs.Range("B3").Copy()
$dbs.Range("B4").PasteSpecial(-4122)
$dbs.Range("C3").Copy()
$dbs.Range("C4").PasteSpecial(-4122)

$dbs.Range("B4").Value = "Entdy >= ,AND Entdy <= ,AND ImportFg = ,AND CustNo >= ,AND CustNo <="
$dbs.Range("A4").Value = "findEntdyImportFg"
$dbs.Range("C4").Value = "CreateDate asc"

# Leave the DBD sheet's last selection parked on B10 ...
$dbd.Select()
$dbd.Range("B10").Select()

# ... and finish with DBS active, selection resting on the next empty row.
$dbs.Select()
$dbs.Range("A5").Select()
